$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header shared-string text updates (preserve surrounding rich-text runs) ---
# A8: "Volume 31   Number  27" -> "...28"
$ws.Range("A8").Characters(21,2).Text = "28"

# C9: "Report Covering the Week  7/1/2024  Through  7/7/2024"
#     -> "...7/8/2024  Through  7/14/2024"
$ws.Range("C9").Characters(27,8).Text = "7/8/2024"
$ws.Range("C9").Characters(46,8).Text = "7/14/2024"

# --- Weekly crime-stat numbers (rows 14-30) ---
$ws.Range("G14").Value = 19
$ws.Range("H14").Value = -57.894736842105
$ws.Range("I14").Value = 58
$ws.Range("J14").Value = 77
$ws.Range("K14").Value = -24.675324675324
$ws.Range("L14").Value = -25.641025641025
$ws.Range("M14").Value = -12.121212121212
$ws.Range("N14").Value = -78.277153558052
$ws.Range("C15").Value = 3
$ws.Range("E15").Value = -62.5
$ws.Range("F15").Value = 28
$ws.Range("G15").Value = 29
$ws.Range("H15").Value = -3.448275862068
$ws.Range("I15").Value = 223
$ws.Range("J15").Value = 212
$ws.Range("K15").Value = 5.188679245283
$ws.Range("L15").Value = 2.293577981651
$ws.Range("M15").Value = 48.666666666666
$ws.Range("N15").Value = -41.005291005291
$ws.Range("C16").Value = 110
$ws.Range("D16").Value = 112
$ws.Range("E16").Value = -1.785714285714
$ws.Range("F16").Value = 430
$ws.Range("G16").Value = 447
$ws.Range("H16").Value = -3.803131991051
$ws.Range("I16").Value = 2614
$ws.Range("J16").Value = 2475
$ws.Range("K16").Value = 5.616161616161
$ws.Range("L16").Value = 0.345489443378
$ws.Range("M16").Value = 15.920177383592
$ws.Range("N16").Value = -69.265138154027
$ws.Range("C17").Value = 181
$ws.Range("D17").Value = 205
$ws.Range("E17").Value = -11.707317073170
$ws.Range("F17").Value = 770
$ws.Range("G17").Value = 783
$ws.Range("H17").Value = -1.660280970625
$ws.Range("I17").Value = 4380
$ws.Range("J17").Value = 4257
$ws.Range("K17").Value = 2.889358703312
$ws.Range("L17").Value = 12.828438948995
$ws.Range("M17").Value = 88.063546586517
$ws.Range("N17").Value = -9.485430874147
$ws.Range("C18").Value = 53
$ws.Range("D18").Value = 65
$ws.Range("E18").Value = -18.461538461538
$ws.Range("F18").Value = 225
$ws.Range("G18").Value = 217
$ws.Range("H18").Value = 3.686635944700
$ws.Range("I18").Value = 1536
$ws.Range("J18").Value = 1603
$ws.Range("K18").Value = -4.179663131628
$ws.Range("L18").Value = -3.578154425612
$ws.Range("M18").Value = -9.593878752207
$ws.Range("N18").Value = -84.559710494571
$ws.Range("C19").Value = 203
$ws.Range("D19").Value = 160
$ws.Range("E19").Value = 26.875
$ws.Range("F19").Value = 704
$ws.Range("G19").Value = 606
$ws.Range("H19").Value = 16.171617161716
$ws.Range("I19").Value = 4749
$ws.Range("J19").Value = 4146
$ws.Range("K19").Value = 14.544138929088
$ws.Range("L19").Value = 12.88328975517
$ws.Range("M19").Value = 100.379746835443
$ws.Range("N19").Value = 23.062969681264
$ws.Range("C20").Value = 85
$ws.Range("D20").Value = 102
$ws.Range("E20").Value = -16.666666666666
$ws.Range("F20").Value = 349
$ws.Range("G20").Value = 376
$ws.Range("H20").Value = -7.180851063829
$ws.Range("I20").Value = 2180
$ws.Range("J20").Value = 2827
$ws.Range("K20").Value = -22.886452069331
$ws.Range("L20").Value = 1.489757914338
$ws.Range("M20").Value = 103.548085901027
$ws.Range("N20").Value = -73.165928114229
$ws.Range("C21").Value = 637
$ws.Range("D21").Value = 656
$ws.Range("E21").Value = -2.896341463414
$ws.Range("F21").Value = 2514
$ws.Range("G21").Value = 2477
$ws.Range("H21").Value = 1.493742430359
$ws.Range("I21").Value = 15740
$ws.Range("J21").Value = 15597
$ws.Range("K21").Value = 0.916842982624
$ws.Range("L21").Value = 6.849501052202
$ws.Range("M21").Value = 58.350100603621
$ws.Range("N21").Value = -56.180400890868
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = -50
$ws.Range("F22").Value = 19
$ws.Range("G22").Value = 16
$ws.Range("H22").Value = 18.75
$ws.Range("I22").Value = 175
$ws.Range("J22").Value = 160
$ws.Range("K22").Value = 9.375
$ws.Range("L22").Value = -10.714285714285
$ws.Range("M22").Value = 1.744186046511
$ws.Range("C23").Value = 21
$ws.Range("D23").Value = 45
$ws.Range("E23").Value = -53.333333333333
$ws.Range("F23").Value = 130
$ws.Range("G23").Value = 153
$ws.Range("H23").Value = -15.032679738562
$ws.Range("I23").Value = 887
$ws.Range("J23").Value = 955
$ws.Range("K23").Value = -7.120418848167
$ws.Range("L23").Value = 2.424942263279
$ws.Range("M23").Value = 59.819819819819
$ws.Range("C24").Value = 340
$ws.Range("D24").Value = 375
$ws.Range("E24").Value = -9.333333333333
$ws.Range("F24").Value = 1154
$ws.Range("G24").Value = 1416
$ws.Range("H24").Value = -18.502824858757
$ws.Range("I24").Value = 8528
$ws.Range("J24").Value = 9496
$ws.Range("K24").Value = -10.193765796124
$ws.Range("L24").Value = -12.970711297071
$ws.Range("M24").Value = 29.133858267716
$ws.Range("C25").Value = 132
$ws.Range("D25").Value = 180
$ws.Range("E25").Value = -26.666666666666
$ws.Range("F25").Value = 449
$ws.Range("G25").Value = 646
$ws.Range("H25").Value = -30.495356037151
$ws.Range("I25").Value = 3396
$ws.Range("J25").Value = 4127
$ws.Range("K25").Value = -17.712624182214
$ws.Range("L25").Value = -33.411764705882
$ws.Range("C26").Value = 235
$ws.Range("D26").Value = 198
$ws.Range("E26").Value = 18.686868686868
$ws.Range("F26").Value = 980
$ws.Range("G26").Value = 845
$ws.Range("H26").Value = 15.976331360946
$ws.Range("I26").Value = 5964
$ws.Range("J26").Value = 5608
$ws.Range("K26").Value = 6.348074179743
$ws.Range("L26").Value = 9.935483870967
$ws.Range("M26").Value = 0.336473755047
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = -40
$ws.Range("F27").Value = 48
$ws.Range("G27").Value = 44
$ws.Range("H27").Value = 9.090909090909
$ws.Range("I27").Value = 349
$ws.Range("J27").Value = 349
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -9.114583333333
$ws.Range("D28").Value = 19
$ws.Range("E28").Value = -5.263157894736
$ws.Range("F28").Value = 89
$ws.Range("G28").Value = 78
$ws.Range("H28").Value = 14.102564102564
$ws.Range("I28").Value = 648
$ws.Range("J28").Value = 563
$ws.Range("K28").Value = 15.097690941385
$ws.Range("L28").Value = 33.884297520661
$ws.Range("C29").Value = 8
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = -20
$ws.Range("F29").Value = 43
$ws.Range("G29").Value = 68
$ws.Range("H29").Value = -36.764705882352
$ws.Range("I29").Value = 208
$ws.Range("J29").Value = 217
$ws.Range("K29").Value = -4.147465437788
$ws.Range("L29").Value = -28.027681660899
$ws.Range("M29").Value = -12.970711297071
$ws.Range("N29").Value = -70.745428973277
$ws.Range("C30").Value = 6
$ws.Range("D30").Value = 7
$ws.Range("E30").Value = -14.285714285714
$ws.Range("F30").Value = 34
$ws.Range("G30").Value = 53
$ws.Range("H30").Value = -35.849056603773
$ws.Range("I30").Value = 169
$ws.Range("J30").Value = 177
$ws.Range("K30").Value = -4.519774011299
$ws.Range("L30").Value = -31.578947368421
$ws.Range("M30").Value = -15.5
$ws.Range("N30").Value = -73.676012461059
# --- Row 31 (Hate Crimes): only the 2-year % change moves ---
$ws.Range("L31").Value = -50

# --- Row 33 (Traffic Fatalities) ---
# C33 flips from a numeric 1 to the literal text "0" (same shared string/style
# already used by D33), so borrow D33's formatting via PasteSpecial after
# writing the text so the style (s=14, not the numeric s=15) matches too.
$ws.Range("C33").Value = "'0"
$ws.Range("D33").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("L33").Value = -31.428571428571
